$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs of "Dessa maneira, haverá um controle maior sobre os
#    dados dos clientes cadastrados." into a single run (Find/Replace over
#    the whole contiguous text re-creates it as one run).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "Dessa maneira, haverá um controle maior sobre os dados dos clientes cadastrados.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Dessa maneira, haverá um controle maior sobre os dados dos clientes cadastrados.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) Append 11 new rows (items 9-19) to the table, re-using the formatting
#    Word copies from the last existing row when Rows.Add() is called.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

$newRows = @(
    @{ Num = "9";  Col2 = "Aplicação de desconto.";   Col3 = "O proprietário terá a opção de aplicar desconto no valor final após finalizar o serviço." },
    @{ Num = "10"; Col2 = "Alteração de dados.";      Col3 = "O proprietário e cliente poderão realizar a alteração dos dados já cadastrados. " },
    @{ Num = "11"; Col2 = "Histórico de serviços.";   Col3 = "O proprietário poderá realizar a verificação dos serviços prestados anteriormente." },
    @{ Num = "12"; Col2 = "Cadastro de horários.";    Col3 = "O proprietário poderá cadastrar os horários disponíveis na agenda." },
    @{ Num = "13"; Col2 = "Cadastro de serviços.";    Col3 = "O proprietário poderá cadastrar os serviços que serão prestados." },
    @{ Num = "14"; Col2 = "Alteração de tema.";       Col3 = "O proprietário e o usuário poderão escolher um tema para aplicação." },
    @{ Num = "15"; Col2 = "Inserir imagem.";          Col3 = "O proprietário poderá inserir imagens de antes, durante e depois do serviço prestado." },
    @{ Num = "16"; Col2 = "Deletar imagem";           Col3 = "O proprietário poderá deletar as imagens já inseridas nos serviços prestados." },
    @{ Num = "17"; Col2 = "Nível de satisfação";      Col3 = "O usuário poderá dar uma nova para o nível do atendimento prestado ao final do serviço." },
    @{ Num = "18"; Col2 = "Comentário";               Col3 = "O proprietário e o usuário poderão deixa um comentário durante a presta ou final do serviço." },
    @{ Num = "19"; Col2 = "Data e hora atual.";       Col3 = "A aplicação irá apresenta a data e hora atual." }
)

foreach ($row in $newRows) {
    $newRow = $t.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $row.Num
    $newRow.Cells.Item(2).Range.Text = $row.Col2
    $newRow.Cells.Item(3).Range.Text = $row.Col3
}

# ---------------------------------------------------------------------------
# 3) The "_GoBack" bookmark used to sit in the paragraph right after the
#    table; it now lives inside row 11 ("Histórico de serviços."), right
#    before the trailing period.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$row11 = $t.Rows.Item($t.Rows.Count - 3)
$col2Cell = $row11.Cells.Item(2)
$bookmarkPos = $col2Cell.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

Write-Output "done"
